$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "in silico" reference columns (OECDQSARTB_* -> inSilico_*/insilico_*)
$ws.Range("J1").Value = "inSilico_AD"
$ws.Range("I1").Value = "insilico_call"

# Shorten the DPRA depletion column headers
$ws.Range("C1").Value = "DPRA_pC"
$ws.Range("D1").Value = "DPRA_pK"

# Match the active selection left in the saved workbook
[void]$ws.Range("D2").Select()
